# Update the workbook to reflect one additional day of data
# (from "through December 23" to "through December 24").
# This updates the sheet name, the header label in B1, and adds/updates
# the December-column counts for each neighborhood row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (also updates xl/workbook.xml sheet name).
$ws.Name = "Through 2021-12-24"

# Update the shared header label for the "through" December 2021 column.
$ws.Range("B1").Value = "December 2021 (through December 24)"

# Update existing counts that increased by one additional day of data.
$ws.Range("BV3").Value = 2
$ws.Range("AX5").Value = 3
$ws.Range("Z6").Value = 5
$ws.Range("N7").Value = 9
$ws.Range("Z7").Value = 6
$ws.Range("B9").Value = 7
$ws.Range("N13").Value = 4
$ws.Range("N15").Value = 8
$ws.Range("BJ17").Value = 2
$ws.Range("B33").Value = 4
$ws.Range("BJ40").Value = 2

# Add new counts for cells that previously had no recorded incidents.
$ws.Range("AL8").Value = 1
$ws.Range("AL12").Value = 1
$ws.Range("BV13").Value = 1
$ws.Range("AX14").Value = 1
$ws.Range("Z26").Value = 1
$ws.Range("BJ44").Value = 1
$ws.Range("N57").Value = 1
$ws.Range("AL61").Value = 1
$ws.Range("B68").Value = 1
$ws.Range("BJ81").Value = 1
$ws.Range("B84").Value = 1
